$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 0.02956397845333333
$ws.Range("R2").Value = 0.26607580608
$ws.Range("S2").Value = 0.006259003216804254
$ws.Range("T2").Value = 0.006259003216804255

# Row 3
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("Q3").Value = 3.502954936099555
$ws.Range("S3").Value = 0.7416121699579786
$ws.Range("T3").Value = 0.7416121699579786

# Row 4
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 1.182964446040889
$ws.Range("R4").Value = 10.646680014368
$ws.Range("S4").Value = 0.2504459365921425
$ws.Range("T4").Value = 0.2504459365921425

# Row 5
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 0.007949018216888889
$ws.Range("R5").Value = 0.07154116395200001
$ws.Range("S5").Value = 0.00168289023307462
$ws.Range("T5").Value = 0.00168289023307462
